$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" numeric-looking string need to be forced
# to stay text (matching the source data, which stores prices as text strings)
# - set a temporary Text number format, assign the value, then restore the
# cell style so no visible formatting/style change is left behind.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '27.584.54'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '1.786.60'
$ws.Range('E3').Value = '  -0.51%  '
Set-TextValue $ws.Range('D4') '1.003'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  -0.03%  '
Set-TextValue $ws.Range('D6') '305.67'
Set-TextValue $ws.Range('D7') '0.4934'
$ws.Range('E7').Value = '  -5.29%  '
Set-TextValue $ws.Range('D8') '0.3827'
$ws.Range('E8').Value = '  +0.42%  '
Set-TextValue $ws.Range('D9') '0.09267'
$ws.Range('E9').Value = '  +16.20%  '
Set-TextValue $ws.Range('D10') '1.085'
$ws.Range('E10').Value = '  -0.93%  '
Set-TextValue $ws.Range('D11') '40.41'
$ws.Range('E11').Value = '  -2.46%  '
Set-TextValue $ws.Range('D12') '1.003'
$ws.Range('E12').Value = '  +0.06%  '
Set-TextValue $ws.Range('D13') '6.203'
$ws.Range('E13').Value = '  -1.42%  '
Set-TextValue $ws.Range('D14') '20.28'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').Value = '1.791.37'
$ws.Range('E15').Value = '  -0.38%  '
Set-TextValue $ws.Range('D16') '7.102'
$ws.Range('E16').Value = '  -2.66%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D17') '0.00001101'
$ws.Range('E17').Value = '  +0.96%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D18') '91.57'
$ws.Range('E18').Value = '  -0.43%  '
Set-TextValue $ws.Range('D19') '0.06524'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('E21').Value = '  -2.25%  '
Set-TextValue $ws.Range('D22') '5.885'
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('D23').Value = '27.653.56'
$ws.Range('E23').Value = '  -1.96%  '
Set-TextValue $ws.Range('D24') '10.86'
$ws.Range('E24').Value = '  -2.51%  '
$ws.Range('E25').Value = '  -2.05%  '
Set-TextValue $ws.Range('D26') '156.62'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').Value = '1.997.72'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('E28').Value = '  -0.86%  '
Set-TextValue $ws.Range('D29') '2.364'
$ws.Range('E29').Value = '  +1.07%  '
Set-TextValue $ws.Range('D30') '125.72'
$ws.Range('E30').Value = '  +2.42%  '
$ws.Range('E31').Value = '  -1.33%  '
Set-TextValue $ws.Range('D32') '1.038'
$ws.Range('E32').Value = '  -1.37%  '
Set-TextValue $ws.Range('D33') '3.606'
$ws.Range('E33').Value = '  -1.82%  '
Set-TextValue $ws.Range('D34') '5.470'
$ws.Range('E34').Value = '  -1.29%  '
Set-TextValue $ws.Range('D35') '0.06757'
$ws.Range('E35').Value = '  -6.69%  '
Set-TextValue $ws.Range('D36') '8.740'
$ws.Range('E36').Value = '  +0.18%  '
Set-TextValue $ws.Range('D37') '0.02280'
$ws.Range('E37').Value = '  -1.75%  '
Set-TextValue $ws.Range('D38') '0.2109'
$ws.Range('E38').Value = '  -1.63%  '
Set-TextValue $ws.Range('D39') '11.26'
$ws.Range('E39').Value = '  -7.35%  '
Set-TextValue $ws.Range('D40') '4.867'
$ws.Range('E40').Value = '  -4.01%  '
Set-TextValue $ws.Range('D41') '0.6065'
$ws.Range('E41').Value = '  -1.47%  '
Set-TextValue $ws.Range('D42') '1.002'
$ws.Range('E42').Value = '  +0.09%  '
Set-TextValue $ws.Range('D43') '1.135'
$ws.Range('E43').Value = '  -2.52%  '
Set-TextValue $ws.Range('D44') '12.97'
$ws.Range('E44').Value = '  -2.15%  '
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D45') '3.652'
$ws.Range('E45').Value = '  -3.14%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D46') '0.5791'
$ws.Range('E46').Value = '  -2.81%  '
Set-TextValue $ws.Range('D47') '1.267'
$ws.Range('E47').Value = '  -6.64%  '
Set-TextValue $ws.Range('D48') '122.29'
$ws.Range('E48').Value = '  -4.63%  '
Set-TextValue $ws.Range('D49') '1.907'
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('E50').Value = '  -6.04%  '
Set-TextValue $ws.Range('D51') '0.06690'
$ws.Range('E51').Value = '  -0.83%  '
